$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert text dates in E2:E7 into real Excel date values with a custom
# date-time number format.
$dates = @(
    @{ Row = 2; Value = 47128 },
    @{ Row = 3; Value = 47578 },
    @{ Row = 4; Value = 44704 },
    @{ Row = 5; Value = 47036 },
    @{ Row = 6; Value = 45544 },
    @{ Row = 7; Value = 45261 }
)

foreach ($d in $dates) {
    $cell = $ws.Cells.Item($d.Row, 5)
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $cell.Value = $d.Value
}
